# Renamed few transcripts. Updated the DataSheet
# Shortens speaker names/tags throughout Sheet1:
#   "Teacher Tag" -> "T Tag"   (header, F1)
#   "Teacher"     -> "T"       (Speaker column D)
#   "Student 1"   -> "S"       (Speaker column D)
#   "Student 2"   -> "SN"      (Speaker column D)
#   "Student 3"   -> "SN"      (Speaker column D)
# and updates the bracketed in-line references inside column E (Sentence)
# to match the shortened names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header tag rename
$ws.Range("F1").Value = "T Tag"

# Map of row -> new speaker (column D)
$speakerMap = @{
    2  = "T";  3  = "T";  4  = "S";  5  = "T";  6  = "S";  7  = "T";  8  = "T";
    9  = "S";  10 = "T"; 11 = "T"; 12 = "T"; 13 = "T"; 14 = "T"; 15 = "T";
    16 = "T"; 17 = "T"; 18 = "T"; 19 = "T"; 20 = "T"; 21 = "T"; 22 = "T";
    23 = "T"; 24 = "T"; 25 = "T"; 26 = "T"; 27 = "T"; 28 = "SN"; 29 = "T";
    30 = "SN"; 31 = "T"; 32 = "T"; 33 = "T"; 34 = "T"; 35 = "SN"; 36 = "T";
    37 = "T"; 38 = "SN"; 39 = "T"; 40 = "T"; 41 = "T"; 42 = "T"; 43 = "T";
    44 = "T"; 45 = "S"; 46 = "T"; 47 = "T"; 48 = "T"; 49 = "SN"; 50 = "T";
    51 = "SN"; 52 = "T"; 53 = "SN"; 54 = "T"; 55 = "T"; 56 = "T"; 57 = "T";
    58 = "T"; 59 = "T"; 60 = "T"; 61 = "T"
}

foreach ($row in $speakerMap.Keys) {
    $ws.Range("D$row").Value = $speakerMap[$row]
}

# Map of row -> new Sentence text (column E) for rows whose bracketed
# student references change along with the name shortening
$sentenceMap = @{
    3  = "[S], would you read it?"
    15 = "So, only two other groups besides [S]."
    19 = "[SN], who's the recorder/reporter?"
    20 = "[SN], come on up."
    21 = "[S], keep yours up there."
    43 = "Um, [SN], hold yours for a minute because that's your person."
    44 = "[S], do you want to rebuttal?"
    48 = "Um, [SN], do you want to rebut him?"
    50 = "Even though [SN] is up here."
    52 = "So, are you saying that your group took [S]'s group number 5 situation and put it with A?"
    55 = "[S], would you take your number 5 just to be visual and would put it over here."
}

foreach ($row in $sentenceMap.Keys) {
    $ws.Range("E$row").Value = $sentenceMap[$row]
}
